$d = $word.ActiveDocument

$d.Content.Find.Execute("260÷6=43, 2", $true, $false, $false, $false, $false, $true, 1, $false, "605÷9=67, 2", 2)
$d.Content.Find.Execute("120÷5=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "624÷9=69, 3", 2)
$d.Content.Find.Execute("433÷2=216, 1", $true, $false, $false, $false, $false, $true, 1, $false, "788÷5=157, 3", 2)
$d.Content.Find.Execute("462÷3=154, 0", $true, $false, $false, $false, $false, $true, 1, $false, "423÷3=141, 0", 2)
$d.Content.Find.Execute("202÷3=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "914÷2=457, 0", 2)
$d.Content.Find.Execute("878÷4=219, 2", $true, $false, $false, $false, $false, $true, 1, $false, "310÷2=155, 0", 2)
$d.Content.Find.Execute("353÷9=39, 2", $true, $false, $false, $false, $false, $true, 1, $false, "543÷6=90, 3", 2)
$d.Content.Find.Execute("286÷2=143, 0", $true, $false, $false, $false, $false, $true, 1, $false, "312÷5=62, 2", 2)
$d.Content.Find.Execute("320÷2=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "988÷8=123, 4", 2)
$d.Content.Find.Execute("179÷5=35, 4", $true, $false, $false, $false, $false, $true, 1, $false, "865÷7=123, 4", 2)
$d.Content.Find.Execute("197÷8=24, 5", $true, $false, $false, $false, $false, $true, 1, $false, "729÷2=364, 1", 2)
$d.Content.Find.Execute("235÷3=78, 1", $true, $false, $false, $false, $false, $true, 1, $false, "868÷8=108, 4", 2)
$d.Content.Find.Execute("783÷8=97, 7", $true, $false, $false, $false, $false, $true, 1, $false, "910÷4=227, 2", 2)
$d.Content.Find.Execute("372÷8=46, 4", $true, $false, $false, $false, $false, $true, 1, $false, "251÷7=35, 6", 2)
$d.Content.Find.Execute("994÷7=142, 0", $true, $false, $false, $false, $false, $true, 1, $false, "625÷2=312, 1", 2)
$d.Content.Find.Execute("374÷6=62, 2", $true, $false, $false, $false, $false, $true, 1, $false, "459÷4=114, 3", 2)
$d.Content.Find.Execute("485÷5=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "396÷7=56, 4", 2)
$d.Content.Find.Execute("947÷6=157, 5", $true, $false, $false, $false, $false, $true, 1, $false, "624÷3=208, 0", 2)
$d.Content.Find.Execute("869÷8=108, 5", $true, $false, $false, $false, $false, $true, 1, $false, "480÷5=96, 0", 2)
$d.Content.Find.Execute("613÷5=122, 3", $true, $false, $false, $false, $false, $true, 1, $false, "103÷8=12, 7", 2)
$d.Content.Find.Execute("705÷7=100, 5", $true, $false, $false, $false, $false, $true, 1, $false, "570÷4=142, 2", 2)
$d.Content.Find.Execute("651÷7=93, 0", $true, $false, $false, $false, $false, $true, 1, $false, "312÷5=62, 2", 2)
$d.Content.Find.Execute("475÷2=237, 1", $true, $false, $false, $false, $false, $true, 1, $false, "901÷9=100, 1", 2)
$d.Content.Find.Execute("162÷7=23, 1", $true, $false, $false, $false, $false, $true, 1, $false, "258÷2=129, 0", 2)
$d.Content.Find.Execute("402÷8=50, 2", $true, $false, $false, $false, $false, $true, 1, $false, "323÷6=53, 5", 2)
